$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1946
$ws.Range("J29").Value = 7500
$ws.Range("L29").Value = 22500
$ws.Range("N29").Value = -23062
$ws.Range("H33").Value = 68.833336
$ws.Range("I33").Value = 73.25
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 73.25
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 155.75
$ws.Range("N33").Value = -518
$ws.Range("H38").Value = 267.25
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 1845.4
$ws.Range("H99").Value = 1715.6364
$ws.Range("I99").Value = 506
$ws.Range("K99").Value = 1518
$ws.Range("M99").Value = -20
$ws.Range("H101").Value = 577.8333
$ws.Range("I101").Value = 804.25
$ws.Range("J101").Value = 125
$ws.Range("K101").Value = 2412.75
$ws.Range("L101").Value = 375
$ws.Range("M101").Value = -790.75
$ws.Range("N101").Value = -3619
$ws.Range("H125").Value = 3123.6667
$ws.Range("I125").Value = 2989.6428
$ws.Range("K125").Value = 26906.7852
$ws.Range("M125").Value = -24446.7852
$ws.Range("H138").Value = 2230.7144
$ws.Range("J138").Value = 2222
$ws.Range("L138").Value = 6666
$ws.Range("N138").Value = -16946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9533.125
$ws.Range("I2").Value = 6838.5
$ws.Range("J2").Value = 12227.75
$ws.Range("K2").Value = 6838.5
$ws.Range("L2").Value = 12227.75
$ws.Range("M2").Value = -6725.5
$ws.Range("N2").Value = -12453.75
$ws.Range("H24").Value = 2508588.5
$ws.Range("J24").Value = 2508588.5
$ws.Range("L24").Value = 2508588.5
$ws.Range("N24").Value = -2509336.5
$ws.Range("H100").Value = 2508588.5
$ws.Range("J100").Value = 2508588.5
$ws.Range("L100").Value = 2508588.5
$ws.Range("N100").Value = -2510752.5
$ws.Range("H102").Value = 4045.56
$ws.Range("I102").Value = 2721.111
$ws.Range("J102").Value = 7451.2856
$ws.Range("K102").Value = 2721.111
$ws.Range("L102").Value = 7451.2856
$ws.Range("M102").Value = -1099.111
$ws.Range("N102").Value = -10695.2856
$ws.Range("H110").Value = 3336.9333
$ws.Range("I110").Value = 2870.9092
$ws.Range("J110").Value = 4618.5
$ws.Range("K110").Value = 2870.9092
$ws.Range("L110").Value = 4618.5
$ws.Range("M110").Value = -825.9092000000001
$ws.Range("N110").Value = -8708.5
$ws.Range("H112").Value = 13666
$ws.Range("J112").Value = 13666
$ws.Range("L112").Value = 13666
$ws.Range("N112").Value = -16620
$ws.Range("H116").Value = 9533.125
$ws.Range("I116").Value = 6838.5
$ws.Range("J116").Value = 12227.75
$ws.Range("K116").Value = 6838.5
$ws.Range("L116").Value = 12227.75
$ws.Range("M116").Value = -4544.5
$ws.Range("N116").Value = -16815.75
$ws.Range("H133").Value = 72761
$ws.Range("J133").Value = 72761
$ws.Range("L133").Value = 72761
$ws.Range("N133").Value = -77821

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9533.125
$ws.Range("I3").Value = 6838.5
$ws.Range("J3").Value = 12227.75
$ws.Range("K3").Value = 6838.5
$ws.Range("L3").Value = 12227.75
$ws.Range("M3").Value = -6724.5
$ws.Range("N3").Value = -12455.75
$ws.Range("H94").Value = 1332.3334
$ws.Range("I94").Value = 1198.8
$ws.Range("K94").Value = 1198.8
$ws.Range("M94").Value = -747.8
$ws.Range("H105").Value = 1366.8334
$ws.Range("I105").Value = 1128.8889
$ws.Range("K105").Value = 1128.8889
$ws.Range("M105").Value = 618.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 952.7143
$ws.Range("I16").Value = 966.6667
$ws.Range("J16").Value = 942.25
$ws.Range("K16").Value = 966.6667
$ws.Range("L16").Value = 942.25
$ws.Range("M16").Value = -679.6667
$ws.Range("N16").Value = -1516.25
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H113").Value = 952.7143
$ws.Range("I113").Value = 966.6667
$ws.Range("J113").Value = 942.25
$ws.Range("K113").Value = 966.6667
$ws.Range("L113").Value = 942.25
$ws.Range("M113").Value = 1203.3333
$ws.Range("N113").Value = -5282.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1864
$ws.Range("J34").Value = 2205.1428
$ws.Range("L34").Value = 6615.428400000001
$ws.Range("N34").Value = -6783.428400000001
$ws.Range("H39").Value = 4773.7144
$ws.Range("J39").Value = 7874.75
$ws.Range("L39").Value = 23624.25
$ws.Range("N39").Value = -24212.25
$ws.Range("H55").Value = 2022
$ws.Range("J55").Value = 4750
$ws.Range("L55").Value = 14250
$ws.Range("N55").Value = -14604
$ws.Range("H60").Value = 1020.25
$ws.Range("I60").Value = 278.92856
$ws.Range("J60").Value = 2750
$ws.Range("K60").Value = 836.78568
$ws.Range("L60").Value = 8250
$ws.Range("M60").Value = -585.78568
$ws.Range("N60").Value = -8752
$ws.Range("H109").Value = 85987.914
$ws.Range("I109").Value = 113328.445
$ws.Range("J109").Value = 3966.3333
$ws.Range("K109").Value = 339985.335
$ws.Range("L109").Value = 11898.9999
$ws.Range("M109").Value = -338945.335
$ws.Range("N109").Value = -13978.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1701.6666
$ws.Range("I80").Value = 1552.5
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1552.5
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -554.5
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 1701.6666
$ws.Range("I83").Value = 1552.5
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 7762.5
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -2770.5
$ws.Range("N83").Value = -19984
$ws.Range("H95").Value = 26666.666
$ws.Range("J95").Value = 26666.666
$ws.Range("L95").Value = 26666.666
$ws.Range("N95").Value = -32158.666
$ws.Range("H102").Value = 4051.3333
$ws.Range("I102").Value = 4051.3333
$ws.Range("K102").Value = 4051.3333
$ws.Range("M102").Value = -2429.3333
$ws.Range("H122").Value = 229624.45
$ws.Range("J122").Value = 3293.6667
$ws.Range("L122").Value = 9881.000100000001
$ws.Range("N122").Value = -14781.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7998
$ws.Range("I7").Value = 6664.6665
$ws.Range("J7").Value = 9998
$ws.Range("K7").Value = 6664.6665
$ws.Range("L7").Value = 9998
$ws.Range("M7").Value = -6552.6665
$ws.Range("N7").Value = -10222
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1280
$ws.Range("H20").Value = 125086.5
$ws.Range("J20").Value = 999999
$ws.Range("L20").Value = 999999
$ws.Range("N20").Value = -1000451
$ws.Range("H61").Value = 6235.375
$ws.Range("H82").Value = 3629.4614
$ws.Range("I82").Value = 2022.875
$ws.Range("K82").Value = 2022.875
$ws.Range("M82").Value = -1661.875
$ws.Range("H85").Value = 3629.4614
$ws.Range("I85").Value = 2022.875
$ws.Range("K85").Value = 2022.875
$ws.Range("M85").Value = -774.875
$ws.Range("H113").Value = 6235.375
$ws.Range("H126").Value = 7998
$ws.Range("I126").Value = 6664.6665
$ws.Range("J126").Value = 9998
$ws.Range("K126").Value = 19993.9995
$ws.Range("L126").Value = 29994
$ws.Range("M126").Value = -17523.9995
$ws.Range("N126").Value = -34934
$ws.Range("H132").Value = 5020.857
$ws.Range("I132").Value = 3594.2
$ws.Range("K132").Value = 10782.6
$ws.Range("M132").Value = -8252.599999999999
$ws.Range("H136").Value = 1500
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2500
$ws.Range("I23").Value = 2500
$ws.Range("K23").Value = 2500
$ws.Range("M23").Value = -2271
$ws.Range("H26").Value = 89999
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H37").Value = 63029
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 63029
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 63029
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -63435
$ws.Range("H113").Value = 935.63635
$ws.Range("I113").Value = 1192.7142
$ws.Range("J113").Value = 485.75
$ws.Range("K113").Value = 3578.1426
$ws.Range("L113").Value = 1457.25
$ws.Range("M113").Value = -1408.1426
$ws.Range("N113").Value = -5797.25
